# Convert a handful of Column B cells from raw numbers (one of them mis-scaled
# due to a missing decimal point) into plain text values that show the
# intended number of decimals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "25.0"

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "37.39"

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "60.0"

$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "60.0"

# Scroll the sheet so row 12 is at the top, and drop the old B2 selection
# (matches the sheetView change: topLeftCell="A12", no <selection> override).
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 12
